$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "day 14" (column O) meal values for each person (rows 3-9)
$ws.Range("O3").Value = 1    # Rakib  - night off
$ws.Range("O4").Value = 2    # Mahfuz
$ws.Range("O5").Value = 1    # Himel  - night off
$ws.Range("O6").Value = 2    # Minhaz
$ws.Range("O7").Value = 2    # Taher
$ws.Range("O8").Value = 2    # Forhad
$ws.Range("O9").Value = 2    # Nayem

# Update the active selection to match the saved cursor position
$ws.Range("O10").Select()

$excel.CalculateFullRebuild()
